$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the title cell (A1, merged A1:E1) to the new, shorter title.
$ws.Range("A1").Value = "Liste tour par tour de distribution des cartes"

# Leave the merged title range selected, as it would be after editing it.
$ws.Range("A1:E1").Select()
